# Update the "Mean flow rate" values (column B, rows 2-9) on Sheet1 with the
# results of the re-run 0.1s model, as described by the commit
# "Modified 0.1s model run complete".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 24.60777723447713
$ws.Range("B3").Value = 1.858667671557108
$ws.Range("B4").Value = 50.622631022042654
$ws.Range("B5").Value = 13.56631717480707
$ws.Range("B6").Value = 3.2040675551705422
$ws.Range("B7").Value = 43.596264765641457
$ws.Range("B8").Value = 17.429177160276236
$ws.Range("B9").Value = 25.137348601619298
